$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 14.734287
$ws.Cells.Item(2, 8).Value = 44.202861
$ws.Cells.Item(2, 9).Value = 0.4000023944294819
$ws.Cells.Item(2, 10).Value = 0.400002394429482
$ws.Cells.Item(2, 13).Value = 4.877755666666666
$ws.Cells.Item(2, 14).Value = 14.633267
$ws.Cells.Item(2, 15).Value = 0.09961167132870688
$ws.Cells.Item(2, 16).Value = 0.09961167132870689
$ws.Cells.Item(2, 17).Value = 71.870251908543
$ws.Cells.Item(2, 18).Value = 646.8322671768869
$ws.Cells.Item(2, 19).Value = 0.03984490704460533
$ws.Cells.Item(2, 20).Value = 0.03984490704460535
$ws.Cells.Item(3, 7).Value = 14.734287
$ws.Cells.Item(3, 8).Value = 44.202861
$ws.Cells.Item(3, 9).Value = 0.4000023944294819
$ws.Cells.Item(3, 10).Value = 0.400002394429482
$ws.Cells.Item(3, 15).Value = 0.1360673938501395
$ws.Cells.Item(3, 16).Value = 0.1360673938501395
$ws.Cells.Item(3, 17).Value = 98.17321346088301
$ws.Cells.Item(3, 18).Value = 883.558921147947
$ws.Cells.Item(3, 19).Value = 0.05442728334383516
$ws.Cells.Item(3, 20).Value = 0.05442728334383519
$ws.Cells.Item(4, 7).Value = 14.734287
$ws.Cells.Item(4, 8).Value = 44.202861
$ws.Cells.Item(4, 9).Value = 0.4000023944294819
$ws.Cells.Item(4, 10).Value = 0.400002394429482
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.246459
$ws.Cells.Item(4, 14).Value = 0.739377
$ws.Cells.Item(4, 15).Value = 0.00503309197542868
$ws.Cells.Item(4, 16).Value = 0.00503309197542868
$ws.Cells.Item(4, 17).Value = 3.631397639733
$ws.Cells.Item(4, 18).Value = 32.682578757597
$ws.Cells.Item(4, 19).Value = 0.002013248841555283
$ws.Cells.Item(4, 20).Value = 0.002013248841555284
$ws.Cells.Item(5, 7).Value = 14.734287
$ws.Cells.Item(5, 8).Value = 44.202861
$ws.Cells.Item(5, 9).Value = 0.4000023944294819
$ws.Cells.Item(5, 10).Value = 0.400002394429482
$ws.Cells.Item(5, 13).Value = 37.01331466666667
$ws.Cells.Item(5, 14).Value = 111.039944
$ws.Cells.Item(5, 15).Value = 0.7558718368280999
$ws.Cells.Item(5, 16).Value = 0.7558718368280999
$ws.Cells.Item(5, 17).Value = 545.3648011199761
$ws.Cells.Item(5, 18).Value = 4908.283210079783
$ws.Cells.Item(5, 19).Value = 0.3023505446130506
$ws.Cells.Item(5, 20).Value = 0.3023505446130507
$ws.Cells.Item(6, 7).Value = 14.734287
$ws.Cells.Item(6, 8).Value = 44.202861
$ws.Cells.Item(6, 9).Value = 0.4000023944294819
$ws.Cells.Item(6, 10).Value = 0.400002394429482
$ws.Cells.Item(6, 13).Value = 0.167274
$ws.Cells.Item(6, 14).Value = 0.501822
$ws.Cells.Item(6, 15).Value = 0.00341600601762507
$ws.Cells.Item(6, 16).Value = 0.00341600601762507
$ws.Cells.Item(6, 17).Value = 2.464663123638
$ws.Cells.Item(6, 18).Value = 22.181968112742
$ws.Cells.Item(6, 19).Value = 0.001366410586435547
$ws.Cells.Item(6, 20).Value = 0.001366410586435547
$ws.Cells.Item(7, 9).Value = 0.3923645715978801
$ws.Cells.Item(7, 10).Value = 0.3923645715978802
$ws.Cells.Item(7, 13).Value = 4.877755666666666
$ws.Cells.Item(7, 14).Value = 14.633267
$ws.Cells.Item(7, 15).Value = 0.09961167132870688
$ws.Cells.Item(7, 16).Value = 0.09961167132870689
$ws.Cells.Item(7, 17).Value = 70.497929496016
$ws.Cells.Item(7, 18).Value = 634.481365464144
$ws.Cells.Item(7, 19).Value = 0.03908409074703691
$ws.Cells.Item(7, 20).Value = 0.03908409074703693
$ws.Cells.Item(8, 9).Value = 0.3923645715978801
$ws.Cells.Item(8, 10).Value = 0.3923645715978802
$ws.Cells.Item(8, 15).Value = 0.1360673938501395
$ws.Cells.Item(8, 16).Value = 0.1360673938501395
$ws.Cells.Item(8, 19).Value = 0.05338802469645001
$ws.Cells.Item(8, 20).Value = 0.05338802469645004
$ws.Cells.Item(9, 9).Value = 0.3923645715978801
$ws.Cells.Item(9, 10).Value = 0.3923645715978802
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.246459
$ws.Cells.Item(9, 14).Value = 0.739377
$ws.Cells.Item(9, 15).Value = 0.00503309197542868
$ws.Cells.Item(9, 16).Value = 0.00503309197542868
$ws.Cells.Item(9, 17).Value = 3.562058125296
$ws.Cells.Item(9, 18).Value = 32.058523127664
$ws.Cells.Item(9, 19).Value = 0.001974806976751802
$ws.Cells.Item(9, 20).Value = 0.001974806976751803
$ws.Cells.Item(10, 9).Value = 0.3923645715978801
$ws.Cells.Item(10, 10).Value = 0.3923645715978802
$ws.Cells.Item(10, 13).Value = 37.01331466666667
$ws.Cells.Item(10, 14).Value = 111.039944
$ws.Cells.Item(10, 15).Value = 0.7558718368280999
$ws.Cells.Item(10, 16).Value = 0.7558718368280999
$ws.Cells.Item(10, 17).Value = 534.951364131712
$ws.Cells.Item(10, 18).Value = 4814.562277185408
$ws.Cells.Item(10, 19).Value = 0.2965773294399602
$ws.Cells.Item(10, 20).Value = 0.2965773294399602
$ws.Cells.Item(11, 9).Value = 0.3923645715978801
$ws.Cells.Item(11, 10).Value = 0.3923645715978802
$ws.Cells.Item(11, 13).Value = 0.167274
$ws.Cells.Item(11, 14).Value = 0.501822
$ws.Cells.Item(11, 15).Value = 0.00341600601762507
$ws.Cells.Item(11, 16).Value = 0.00341600601762507
$ws.Cells.Item(11, 17).Value = 2.417601754656
$ws.Cells.Item(11, 18).Value = 21.758415791904
$ws.Cells.Item(11, 19).Value = 0.001340319737681241
$ws.Cells.Item(11, 20).Value = 0.001340319737681242
$ws.Cells.Item(12, 7).Value = 1.259379333333333
$ws.Cells.Item(12, 8).Value = 3.778138
$ws.Cells.Item(12, 9).Value = 0.03418928576783783
$ws.Cells.Item(12, 10).Value = 0.03418928576783784
$ws.Cells.Item(12, 13).Value = 4.877755666666666
$ws.Cells.Item(12, 14).Value = 14.633267
$ws.Cells.Item(12, 15).Value = 0.09961167132870688
$ws.Cells.Item(12, 16).Value = 0.09961167132870689
$ws.Cells.Item(12, 17).Value = 6.142944679649555
$ws.Cells.Item(12, 18).Value = 55.286502116846
$ws.Cells.Item(12, 19).Value = 0.003405651896869098
$ws.Cells.Item(12, 20).Value = 0.003405651896869099
$ws.Cells.Item(13, 7).Value = 1.259379333333333
$ws.Cells.Item(13, 8).Value = 3.778138
$ws.Cells.Item(13, 9).Value = 0.03418928576783783
$ws.Cells.Item(13, 10).Value = 0.03418928576783784
$ws.Cells.Item(13, 15).Value = 0.1360673938501395
$ws.Cells.Item(13, 16).Value = 0.1360673938501395
$ws.Cells.Item(13, 17).Value = 8.391129894480667
$ws.Cells.Item(13, 18).Value = 75.520169050326
$ws.Cells.Item(13, 19).Value = 0.004652047012027359
$ws.Cells.Item(13, 20).Value = 0.004652047012027361
$ws.Cells.Item(14, 7).Value = 1.259379333333333
$ws.Cells.Item(14, 8).Value = 3.778138
$ws.Cells.Item(14, 9).Value = 0.03418928576783783
$ws.Cells.Item(14, 10).Value = 0.03418928576783784
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.246459
$ws.Cells.Item(14, 14).Value = 0.739377
$ws.Cells.Item(14, 15).Value = 0.00503309197542868
$ws.Cells.Item(14, 16).Value = 0.00503309197542868
$ws.Cells.Item(14, 17).Value = 0.310385371114
$ws.Cells.Item(14, 18).Value = 2.793468340026
$ws.Cells.Item(14, 19).Value = 0.0001720778198437426
$ws.Cells.Item(14, 20).Value = 0.0001720778198437426
$ws.Cells.Item(15, 7).Value = 1.259379333333333
$ws.Cells.Item(15, 8).Value = 3.778138
$ws.Cells.Item(15, 9).Value = 0.03418928576783783
$ws.Cells.Item(15, 10).Value = 0.03418928576783784
$ws.Cells.Item(15, 13).Value = 37.01331466666667
$ws.Cells.Item(15, 14).Value = 111.039944
$ws.Cells.Item(15, 15).Value = 0.7558718368280999
$ws.Cells.Item(15, 16).Value = 0.7558718368280999
$ws.Cells.Item(15, 17).Value = 46.61380354936356
$ws.Cells.Item(15, 18).Value = 419.524231944272
$ws.Cells.Item(15, 19).Value = 0.0258427182331764
$ws.Cells.Item(15, 20).Value = 0.0258427182331764
$ws.Cells.Item(16, 7).Value = 1.259379333333333
$ws.Cells.Item(16, 8).Value = 3.778138
$ws.Cells.Item(16, 9).Value = 0.03418928576783783
$ws.Cells.Item(16, 10).Value = 0.03418928576783784
$ws.Cells.Item(16, 13).Value = 0.167274
$ws.Cells.Item(16, 14).Value = 0.501822
$ws.Cells.Item(16, 15).Value = 0.00341600601762507
$ws.Cells.Item(16, 16).Value = 0.00341600601762507
$ws.Cells.Item(16, 17).Value = 0.210661418604
$ws.Cells.Item(16, 18).Value = 1.895952767436
$ws.Cells.Item(16, 19).Value = 0.0001167908059212372
$ws.Cells.Item(16, 20).Value = 0.0001167908059212372
$ws.Cells.Item(17, 7).Value = 4.524255666666667
$ws.Cells.Item(17, 8).Value = 13.572767
$ws.Cells.Item(17, 9).Value = 0.1228232556945456
$ws.Cells.Item(17, 10).Value = 0.1228232556945456
$ws.Cells.Item(17, 13).Value = 4.877755666666666
$ws.Cells.Item(17, 14).Value = 14.633267
$ws.Cells.Item(17, 15).Value = 0.09961167132870688
$ws.Cells.Item(17, 16).Value = 0.09961167132870689
$ws.Cells.Item(17, 17).Value = 22.06821371553211
$ws.Cells.Item(17, 18).Value = 198.613923439789
$ws.Cells.Item(17, 19).Value = 0.0122346297777668
$ws.Cells.Item(17, 20).Value = 0.0122346297777668
$ws.Cells.Item(18, 7).Value = 4.524255666666667
$ws.Cells.Item(18, 8).Value = 13.572767
$ws.Cells.Item(18, 9).Value = 0.1228232556945456
$ws.Cells.Item(18, 10).Value = 0.1228232556945456
$ws.Cells.Item(18, 15).Value = 0.1360673938501395
$ws.Cells.Item(18, 16).Value = 0.1360673938501395
$ws.Cells.Item(18, 17).Value = 30.14470379973433
$ws.Cells.Item(18, 18).Value = 271.302334197609
$ws.Cells.Item(18, 19).Value = 0.01671224030654612
$ws.Cells.Item(18, 20).Value = 0.01671224030654613
$ws.Cells.Item(19, 7).Value = 4.524255666666667
$ws.Cells.Item(19, 8).Value = 13.572767
$ws.Cells.Item(19, 9).Value = 0.1228232556945456
$ws.Cells.Item(19, 10).Value = 0.1228232556945456
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.246459
$ws.Cells.Item(19, 14).Value = 0.739377
$ws.Cells.Item(19, 15).Value = 0.00503309197542868
$ws.Cells.Item(19, 16).Value = 0.00503309197542868
$ws.Cells.Item(19, 17).Value = 1.115043527351
$ws.Cells.Item(19, 18).Value = 10.035391746159
$ws.Cells.Item(19, 19).Value = 0.0006181807426322422
$ws.Cells.Item(19, 20).Value = 0.0006181807426322423
$ws.Cells.Item(20, 7).Value = 4.524255666666667
$ws.Cells.Item(20, 8).Value = 13.572767
$ws.Cells.Item(20, 9).Value = 0.1228232556945456
$ws.Cells.Item(20, 10).Value = 0.1228232556945456
$ws.Cells.Item(20, 13).Value = 37.01331466666667
$ws.Cells.Item(20, 14).Value = 111.039944
$ws.Cells.Item(20, 15).Value = 0.7558718368280999
$ws.Cells.Item(20, 16).Value = 0.7558718368280999
$ws.Cells.Item(20, 17).Value = 167.4576986227831
$ws.Cells.Item(20, 18).Value = 1507.119287605048
$ws.Cells.Item(20, 19).Value = 0.09283863988704355
$ws.Cells.Item(20, 20).Value = 0.09283863988704356
$ws.Cells.Item(21, 7).Value = 4.524255666666667
$ws.Cells.Item(21, 8).Value = 13.572767
$ws.Cells.Item(21, 9).Value = 0.1228232556945456
$ws.Cells.Item(21, 10).Value = 0.1228232556945456
$ws.Cells.Item(21, 13).Value = 0.167274
$ws.Cells.Item(21, 14).Value = 0.501822
$ws.Cells.Item(21, 15).Value = 0.00341600601762507
$ws.Cells.Item(21, 16).Value = 0.00341600601762507
$ws.Cells.Item(21, 17).Value = 0.7567903423860001
$ws.Cells.Item(21, 18).Value = 6.811113081474001
$ws.Cells.Item(21, 19).Value = 0.0004195649805568703
$ws.Cells.Item(21, 20).Value = 0.0004195649805568704
$ws.Cells.Item(22, 7).Value = 1.864631
$ws.Cells.Item(22, 8).Value = 5.593893
$ws.Cells.Item(22, 9).Value = 0.05062049251025444
$ws.Cells.Item(22, 10).Value = 0.05062049251025445
$ws.Cells.Item(22, 13).Value = 4.877755666666666
$ws.Cells.Item(22, 14).Value = 14.633267
$ws.Cells.Item(22, 15).Value = 0.09961167132870688
$ws.Cells.Item(22, 16).Value = 0.09961167132870689
$ws.Cells.Item(22, 17).Value = 9.095214426492332
$ws.Cells.Item(22, 18).Value = 81.856929838431
$ws.Cells.Item(22, 19).Value = 0.005042391862428734
$ws.Cells.Item(22, 20).Value = 0.005042391862428735
$ws.Cells.Item(23, 7).Value = 1.864631
$ws.Cells.Item(23, 8).Value = 5.593893
$ws.Cells.Item(23, 9).Value = 0.05062049251025444
$ws.Cells.Item(23, 10).Value = 0.05062049251025445
$ws.Cells.Item(23, 15).Value = 0.1360673938501395
$ws.Cells.Item(23, 16).Value = 0.1360673938501395
$ws.Cells.Item(23, 17).Value = 12.423866671579
$ws.Cells.Item(23, 18).Value = 111.814800044211
$ws.Cells.Item(23, 19).Value = 0.006887798491280827
$ws.Cells.Item(23, 20).Value = 0.00688779849128083
$ws.Cells.Item(24, 7).Value = 1.864631
$ws.Cells.Item(24, 8).Value = 5.593893
$ws.Cells.Item(24, 9).Value = 0.05062049251025444
$ws.Cells.Item(24, 10).Value = 0.05062049251025445
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 12).Value = 1
$ws.Cells.Item(24, 13).Value = 0.246459
$ws.Cells.Item(24, 14).Value = 0.739377
$ws.Cells.Item(24, 15).Value = 0.00503309197542868
$ws.Cells.Item(24, 16).Value = 0.00503309197542868
$ws.Cells.Item(24, 17).Value = 0.4595550916289999
$ws.Cells.Item(24, 18).Value = 4.135995824660999
$ws.Cells.Item(24, 19).Value = 0.0002547775946456092
$ws.Cells.Item(24, 20).Value = 0.0002547775946456092
$ws.Cells.Item(25, 7).Value = 1.864631
$ws.Cells.Item(25, 8).Value = 5.593893
$ws.Cells.Item(25, 9).Value = 0.05062049251025444
$ws.Cells.Item(25, 10).Value = 0.05062049251025445
$ws.Cells.Item(25, 13).Value = 37.01331466666667
$ws.Cells.Item(25, 14).Value = 111.039944
$ws.Cells.Item(25, 15).Value = 0.7558718368280999
$ws.Cells.Item(25, 16).Value = 0.7558718368280999
$ws.Cells.Item(25, 17).Value = 69.01617394022134
$ws.Cells.Item(25, 18).Value = 621.1455654619919
$ws.Cells.Item(25, 19).Value = 0.0382626046548691
$ws.Cells.Item(25, 20).Value = 0.0382626046548691
$ws.Cells.Item(26, 7).Value = 1.864631
$ws.Cells.Item(26, 8).Value = 5.593893
$ws.Cells.Item(26, 9).Value = 0.05062049251025444
$ws.Cells.Item(26, 10).Value = 0.05062049251025445
$ws.Cells.Item(26, 13).Value = 0.167274
$ws.Cells.Item(26, 14).Value = 0.501822
$ws.Cells.Item(26, 15).Value = 0.00341600601762507
$ws.Cells.Item(26, 16).Value = 0.00341600601762507
$ws.Cells.Item(26, 17).Value = 0.311904285894
$ws.Cells.Item(26, 18).Value = 2.807138573046
$ws.Cells.Item(26, 19).Value = 0.000172919907030174
$ws.Cells.Item(26, 20).Value = 0.000172919907030174
